# "contig stats and figures"
# Adds a Stress column to the existing TAXONOMY block, and appends three
# new blocks of Bray-Curtis/Jaccard distance stats (COG, ARG, CAZymes),
# mirroring the layout of the pre-existing TAXONOMY block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TAXONOMY block: add new "Stress" column (D) ---
$ws.Range('D4').Value = 'Stress'
$ws.Range('D5').Value = 0.124
$ws.Range('D6').Value = 6.5

# --- COG block (rows 18-33) ---
$ws.Range('A18').Value = 'COG'
$ws.Range('A19').Value = 'Bray-Curtis'
$ws.Range('A20').Value = 'Fecal'
$ws.Range('I20').Value = 'Oral'
$ws.Range('B21').Value = 'Ethnicity'
$ws.Range('D21').Value = 'Stress'
$ws.Range('F21').Value = 'Study period'
$ws.Range('J21').Value = 'Ethnicity'
$ws.Range('N21').Value = 'Study period'
$ws.Range('A22').Value = 'All'
$ws.Range('B22').Value = 0.018
$ws.Range('C22').Value = 0.146
$ws.Range('D22').Value = 0.075
$ws.Range('E22').Value = 'All'
$ws.Range('F22').Value = 0.021
$ws.Range('G22').Value = 0.331
$ws.Range('I22').Value = 'All'
$ws.Range('J22').Value = 0.038
$ws.Range('K22').Value = 0.269
$ws.Range('M22').Value = 'All'
$ws.Range('N22').Value = 0.019
$ws.Range('O22').Value = 0.462
$ws.Range('A23').Value = 'Before'
$ws.Range('B23').Value = 0.118
$ws.Range('C23').Value = 0.143
$ws.Range('E23').Value = 'Black'
$ws.Range('F23').Value = 0.011
$ws.Range('G23').Value = 0.951
$ws.Range('I23').Value = 'Before'
$ws.Range('J23').Value = 0.072
$ws.Range('K23').Value = 0.289
$ws.Range('M23').Value = 'Black'
$ws.Range('N23').Value = 0.057
$ws.Range('O23').Value = 0.417
$ws.Range('A24').Value = 'During'
$ws.Range('B24').Value = 0.016
$ws.Range('C24').Value = 0.354
$ws.Range('E24').Value = 'White'
$ws.Range('F24').Value = 0.048
$ws.Range('G24').Value = 0.202
$ws.Range('I24').Value = 'After'
$ws.Range('J24').Value = 0.079
$ws.Range('K24').Value = 0.271
$ws.Range('M24').Value = 'White'
$ws.Range('N24').Value = 0.017
$ws.Range('O24').Value = 0.638
$ws.Range('A25').Value = 'After'
$ws.Range('B25').Value = 0.032
$ws.Range('C25').Value = 0.356
$ws.Range('I26').Value = 'Jaccard'
$ws.Range('A27').Value = 'Jaccard'
$ws.Range('I27').Value = 'Fecal'
$ws.Range('A28').Value = 'Fecal'
$ws.Range('B29').Value = 'Ethnicity'
$ws.Range('F29').Value = 'Study period'
$ws.Range('J29').Value = 'Ethnicity'
$ws.Range('N29').Value = 'Study period'
$ws.Range('A30').Value = 'All'
$ws.Range('B30').Value = 0.022
$ws.Range('C30').Value = 0.096
$ws.Range('E30').Value = 'All'
$ws.Range('F30').Value = 0.019
$ws.Range('G30').Value = 0.374
$ws.Range('I30').Value = 'All'
$ws.Range('J30').Value = 0.072
$ws.Range('K30').Value = 0.086
$ws.Range('M30').Value = 'All'
$ws.Range('N30').Value = 0.019
$ws.Range('O30').Value = 0.545
$ws.Range('A31').Value = 'Before'
$ws.Range('B31').Value = 0.118
$ws.Range('C31').Value = 0.139
$ws.Range('E31').Value = 'Black'
$ws.Range('F31').Value = 0.012
$ws.Range('G31').Value = 0.972
$ws.Range('I31').Value = 'Before'
$ws.Range('J31').Value = 0.072
$ws.Range('K31').Value = 0.294
$ws.Range('M31').Value = 'Black'
$ws.Range('N31').Value = 0.051
$ws.Range('O31').Value = 0.47
$ws.Range('A32').Value = 'During'
$ws.Range('B32').Value = 0.016
$ws.Range('C32').Value = 0.361
$ws.Range('E32').Value = 'White'
$ws.Range('F32').Value = 0.043
$ws.Range('G32').Value = 0.237
$ws.Range('I32').Value = 'After'
$ws.Range('J32').Value = 0.079
$ws.Range('K32').Value = 0.264
$ws.Range('M32').Value = 'White'
$ws.Range('N32').Value = 0.018
$ws.Range('O32').Value = 0.721
$ws.Range('A33').Value = 'After'
$ws.Range('B33').Value = 0.032
$ws.Range('C33').Value = 0.333

# --- ARG block (rows 36-51) ---
$ws.Range('A36').Value = 'ARG'
$ws.Range('A37').Value = 'Bray-Curtis'
$ws.Range('A38').Value = 'Fecal'
$ws.Range('I38').Value = 'Oral'
$ws.Range('B39').Value = 'Ethnicity'
$ws.Range('D39').Value = 'Stress'
$ws.Range('F39').Value = 'Study period'
$ws.Range('J39').Value = 'Ethnicity'
$ws.Range('N39').Value = 'Study period'
$ws.Range('A40').Value = 'All'
$ws.Range('B40').Value = 0.006
$ws.Range('C40').Value = 0.536
$ws.Range('D40').Value = 0.06
$ws.Range('E40').Value = 'All'
$ws.Range('I40').Value = 'All'
$ws.Range('J40').Value = 0.025
$ws.Range('K40').Value = 0.391
$ws.Range('M40').Value = 'All'
$ws.Range('A41').Value = 'Before'
$ws.Range('B41').Value = 0.025
$ws.Range('C41').Value = 0.683
$ws.Range('E41').Value = 'Black'
$ws.Range('I41').Value = 'Before'
$ws.Range('J41').Value = 0.048
$ws.Range('K41').Value = 0.472
$ws.Range('M41').Value = 'Black'
$ws.Range('A42').Value = 'During'
$ws.Range('B42').Value = 0.008
$ws.Range('C42').Value = 0.674
$ws.Range('E42').Value = 'White'
$ws.Range('I42').Value = 'After'
$ws.Range('J42').Value = 0.067
$ws.Range('K42').Value = 0.322
$ws.Range('M42').Value = 'White'
$ws.Range('A43').Value = 'After'
$ws.Range('B43').Value = 0.012
$ws.Range('C43').Value = 0.822
$ws.Range('I44').Value = 'Jaccard'
$ws.Range('A45').Value = 'Jaccard'
$ws.Range('I45').Value = 'Fecal'
$ws.Range('A46').Value = 'Fecal'
$ws.Range('B47').Value = 'Ethnicity'
$ws.Range('F47').Value = 'Study period'
$ws.Range('J47').Value = 'Ethnicity'
$ws.Range('N47').Value = 'Study period'
$ws.Range('A48').Value = 'All'
$ws.Range('B48').Value = 0.006
$ws.Range('C48').Value = 0.584
$ws.Range('E48').Value = 'All'
$ws.Range('I48').Value = 'All'
$ws.Range('J48').Value = 0.054
$ws.Range('K48').Value = 0.148
$ws.Range('M48').Value = 'All'
$ws.Range('A49').Value = 'Before'
$ws.Range('B49').Value = 0.025
$ws.Range('C49').Value = 0.683
$ws.Range('E49').Value = 'Black'
$ws.Range('I49').Value = 'Before'
$ws.Range('J49').Value = 0.048
$ws.Range('K49').Value = 0.473
$ws.Range('M49').Value = 'Black'
$ws.Range('A50').Value = 'During'
$ws.Range('B50').Value = 0.008
$ws.Range('C50').Value = 0.677
$ws.Range('E50').Value = 'White'
$ws.Range('I50').Value = 'After'
$ws.Range('J50').Value = 0.067
$ws.Range('K50').Value = 0.345
$ws.Range('M50').Value = 'White'
$ws.Range('A51').Value = 'After'
$ws.Range('B51').Value = 0.012
$ws.Range('C51').Value = 0.825

# --- CAZymes block (rows 53-57) ---
$ws.Range('A53').Value = 'CAZymes'
$ws.Range('A54').Value = 'Bray-Curtis'
$ws.Range('A55').Value = 'Fecal'
$ws.Range('I55').Value = 'Oral'
$ws.Range('B56').Value = 'Ethnicity'
$ws.Range('D56').Value = 'Stress'
$ws.Range('F56').Value = 'Study period'
$ws.Range('J56').Value = 'Ethnicity'
$ws.Range('L56').Value = 'Stress'
$ws.Range('N56').Value = 'Study period'
$ws.Range('A57').Value = 'All'
$ws.Range('B57').Value = 0.017
$ws.Range('C57').Value = 0.154
$ws.Range('D57').Value = 0.049
$ws.Range('E57').Value = 'All'
$ws.Range('I57').Value = 'All'
$ws.Range('J57').Value = 0.099
$ws.Range('K57').Value = 0.046
$ws.Range('L57').Value = 0.072
$ws.Range('M57').Value = 'All'

# Leave the view scrolled to / selecting the last-edited cell, as in the
# saved workbook.
$ws.Range('L57').Select()
